# Updates cryptos list values (Price and Volume(1h) columns) to reflect
# the latest scrape, per commit "Updated cryptos list on Tue Mar 21 20:32:31 UTC 2023 with GitHub Actions".
# Also swaps the TheSandbox/Hedera rows (37 and 39) to reflect updated ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value2 = "28.248.76"
$dCell.Style = "Normal"
$ws.Range("E2").Value2 = "  +0.95%  "

# Row 3
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.798.52"
$dCell.Style = "Normal"
$ws.Range("E3").Value2 = "  +2.08%  "

# Row 4
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.002"
$dCell.Style = "Normal"
$ws.Range("E4").Value2 = "  -0.04%  "

# Row 5
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value2 = "339.00"
$dCell.Style = "Normal"
$ws.Range("E5").Value2 = "  +0.67%  "

# Row 6
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.000"
$dCell.Style = "Normal"
$ws.Range("E6").Value2 = "  +0.20%  "

# Row 7
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.4800"
$dCell.Style = "Normal"
$ws.Range("E7").Value2 = "  +27.16%  "

# Row 8
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.3690"
$dCell.Style = "Normal"
$ws.Range("E8").Value2 = "  +9.94%  "

# Row 9
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value2 = "45.49"
$dCell.Style = "Normal"
$ws.Range("E9").Value2 = "  -0.39%  "

# Row 10
$ws.Range("E10").Value2 = "  +6.81%  "

# Row 11
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.145"
$dCell.Style = "Normal"
$ws.Range("E11").Value2 = "  +2.05%  "

# Row 12
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value2 = "22.59"
$dCell.Style = "Normal"
$ws.Range("E12").Value2 = "  +1.08%  "

# Row 13
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.003"
$dCell.Style = "Normal"
$ws.Range("E13").Value2 = "  +0.21%  "

# Row 14
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value2 = "6.301"
$dCell.Style = "Normal"
$ws.Range("E14").Value2 = "  +1.65%  "

# Row 15
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value2 = "7.317"
$dCell.Style = "Normal"
$ws.Range("E15").Value2 = "  +1.86%  "

# Row 16
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.798.40"
$dCell.Style = "Normal"
$ws.Range("E16").Value2 = "  +2.38%  "

# Row 17
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.00001097"
$dCell.Style = "Normal"
$ws.Range("E17").Value2 = "  +4.08%  "

# Row 18
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.06708"
$dCell.Style = "Normal"
$ws.Range("E18").Value2 = "  +1.94%  "

# Row 19
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value2 = "82.13"
$dCell.Style = "Normal"
$ws.Range("E19").Value2 = "  +2.04%  "

# Row 20
$ws.Range("E20").Value2 = "  +0.03%  "

# Row 21
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value2 = "17.39"
$dCell.Style = "Normal"
$ws.Range("E21").Value2 = "  +2.29%  "

# Row 22
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value2 = "6.407"
$dCell.Style = "Normal"
$ws.Range("E22").Value2 = "  +1.93%  "

# Row 23
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value2 = "28.231.63"
$dCell.Style = "Normal"
$ws.Range("E23").Value2 = "  +0.85%  "

# Row 24
$ws.Range("E24").Value2 = "  +2.60%  "

# Row 25
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value2 = "2.399"
$dCell.Style = "Normal"
$ws.Range("E25").Value2 = "  +1.45%  "

# Row 26
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value2 = "20.78"
$dCell.Style = "Normal"
$ws.Range("E26").Value2 = "  +4.62%  "

# Row 27
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value2 = "2.410"
$dCell.Style = "Normal"
$ws.Range("E27").Value2 = "  +2.94%  "

# Row 28
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value2 = "151.01"
$dCell.Style = "Normal"
$ws.Range("E28").Value2 = "  -1.16%  "

# Row 29
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value2 = "2.003.36"
$dCell.Style = "Normal"
$ws.Range("E29").Value2 = "  +2.31%  "

# Row 30
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value2 = "134.19"
$dCell.Style = "Normal"
$ws.Range("E30").Value2 = "  +1.63%  "

# Row 31
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.269"
$dCell.Style = "Normal"
$ws.Range("E31").Value2 = "  +1.01%  "

# Row 32
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value2 = "4.037"
$dCell.Style = "Normal"
$ws.Range("E32").Value2 = "  +0.55%  "

# Row 33
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.09669"
$dCell.Style = "Normal"
$ws.Range("E33").Value2 = "  +10.05%  "

# Row 34
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value2 = "5.940"
$dCell.Style = "Normal"
$ws.Range("E34").Value2 = "  +2.09%  "

# Row 35
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.02373"
$dCell.Style = "Normal"
$ws.Range("E35").Value2 = "  +1.20%  "

# Row 36
$ws.Range("E36").Value2 = "  -0.96%  "

# Row 37
$ws.Range("B37").Value2 = "Hedera"
$ws.Range("C37").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.06321"
$dCell.Style = "Normal"
$ws.Range("E37").Value2 = "  +2.28%  "

# Row 38
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.2191"
$dCell.Style = "Normal"
$ws.Range("E38").Value2 = "  +3.66%  "

# Row 39
$ws.Range("B39").Value2 = "TheSandbox"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.6677"
$dCell.Style = "Normal"
$ws.Range("E39").Value2 = "  +1.10%  "

# Row 40
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value2 = "5.231"
$dCell.Style = "Normal"
$ws.Range("E40").Value2 = "  +1.29%  "

# Row 41
$ws.Range("E41").Value2 = "  +2.77%  "

# Row 42
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.221"
$dCell.Style = "Normal"
$ws.Range("E42").Value2 = "  +0.87%  "

# Row 43
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value2 = "8.115"
$dCell.Style = "Normal"
$ws.Range("E43").Value2 = "  +1.13%  "

# Row 44
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value2 = "14.16"
$dCell.Style = "Normal"
$ws.Range("E44").Value2 = "  +3.04%  "

# Row 45
$ws.Range("E45").Value2 = "  +0.16%  "

# Row 46
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.6157"
$dCell.Style = "Normal"
$ws.Range("E46").Value2 = "  +1.62%  "

# Row 47
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value2 = "3.869"
$dCell.Style = "Normal"
$ws.Range("E47").Value2 = "  +1.18%  "

# Row 48
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value2 = "128.89"
$dCell.Style = "Normal"
$ws.Range("E48").Value2 = "  -1.06%  "

# Row 49
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value2 = "2.048"
$dCell.Style = "Normal"
$ws.Range("E49").Value2 = "  +1.69%  "

# Row 50
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value2 = "1.171"
$dCell.Style = "Normal"
$ws.Range("E50").Value2 = "  -1.16%  "

# Row 51
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value2 = "0.07099"
$dCell.Style = "Normal"
$ws.Range("E51").Value2 = "  -0.91%  "

